$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared StartDate / EndDate strings used by rows 2-9
# (20/02/2019 -> 10/02/2019, 22/02/2019 -> 12/02/2019)
$ws.Range("C2:C9").Value = "10/02/2019"
$ws.Range("D2:D9").Value = "12/02/2019"

# Update Kids / AgeKid1 counts for rows 6-9 from 1/5 to 0/0.
# The cells are formatted as text ("@"), so temporarily switch to a
# general number format to make sure the values are written as numbers,
# then restore the original text format.
$ws.Range("F6:G9").NumberFormat = "General"
$ws.Range("F6:G9").Value = 0
$ws.Range("F6:G9").NumberFormat = "@"

# Reset the sheet view/selection back to A1 (was topLeftCell="B1",
# activeCell="M1", sqref="M1:O17")
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1").Select()
